$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# "About" sheet: update version text in A2 and citation text in A6
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "Version: $newVersion"
$wsAbout.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Integra Underground Mine, Australia, M4354, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# "Boundaries and methane sources" sheet: update build_version column (S) for every data row
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
$lastRow = $wsData.Cells.Item($wsData.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
